# Updated NATMI TPM recompute for the Wnt16-Fzd2 LR pair (YoungD7).
#
# The ligand/receptor expression statistics (columns G:T) were regenerated
# from the new TPM matrix for the three remaining sending clusters (ECs,
# FAPs, MuSCs). The "Resolving-Mac" sending-cluster block (the former rows
# 14-17, where Resolving-Mac itself was the "Sending cluster") no longer
# exists in the recomputed output, so it is deleted below. "Resolving-Mac"
# is still a valid "Target cluster" value and remains in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# New column G:T values per data row (rows 2-13), recomputed with the new TPM data.
$newValues = @{
    2  = @(0.702341, 2.107023, 0.3978241315558134, 0.3978241315558134, 3.0, 1.0, 0.1503136666666667, 0.450941, 0.008029526741163598, 0.008029526741163598, 0.1055714509603333, 0.950143058643, 0.003194339502607588, 0.003194339502607588)
    3  = @(0.702341, 2.107023, 0.3978241315558134, 0.3978241315558134, 3.0, 1.0, 15.70489833333333, 47.114695, 0.8389317081486641, 0.8389317081486641, 11.03019400033167, 99.27174600298498, 0.3337472782288774, 0.3337472782288774)
    4  = @(0.702341, 2.107023, 0.3978241315558134, 0.3978241315558134, 3.0, 1.0, 2.758466666666667, 8.275399999999999, 0.1473530807662759, 0.1473530807662759, 1.937384237133333, 17.4364581342, 0.05862061138791733, 0.05862061138791733)
    5  = @(0.702341, 2.107023, 0.3978241315558134, 0.3978241315558134, 3.0, 1.0, 0.1064366666666667, 0.31931, 0.005685684343896315, 0.005685684343896314, 0.07475483490333333, 0.6727935141299999, 0.002261902436411036, 0.002261902436411035)
    6  = @(0.955321, 2.865963, 0.5411185552061337, 0.5411185552061337, 3.0, 1.0, 0.1503136666666667, 0.450941, 0.008029526741163598, 0.008029526741163598, 0.1435978023536667, 1.292380221183, 0.004344925909167461, 0.004344925909167461)
    7  = @(0.955321, 2.865963, 0.5411185552061337, 0.5411185552061337, 3.0, 1.0, 15.70489833333333, 47.114695, 0.8389317081486641, 0.8389317081486641, 15.00321918069833, 135.028972626285, 0.453961513830019, 0.453961513830019)
    8  = @(0.955321, 2.865963, 0.5411185552061337, 0.5411185552061337, 3.0, 1.0, 2.758466666666667, 8.275399999999999, 0.1473530807662759, 0.1473530807662759, 2.635221134466667, 23.7169902102, 0.07973548616941994, 0.07973548616941994)
    9  = @(0.955321, 2.865963, 0.5411185552061337, 0.5411185552061337, 3.0, 1.0, 0.1064366666666667, 0.31931, 0.005685684343896315, 0.005685684343896314, 0.1016811828366667, 0.9151306455299999, 0.003076629297527308, 0.003076629297527308)
    10 = @(0.107794, 0.323382, 0.06105731323805295, 0.06105731323805295, 3.0, 1.0, 0.1503136666666667, 0.450941, 0.008029526741163598, 0.008029526741163598, 0.01620291138466667, 0.145826202462, 0.0004902613293885483, 0.0004902613293885483)
    11 = @(0.107794, 0.323382, 0.06105731323805295, 0.06105731323805295, 3.0, 1.0, 15.70489833333333, 47.114695, 0.8389317081486641, 0.8389317081486641, 1.692893810943333, 15.23604429849, 0.0512229160897678, 0.05122291608976781)
    12 = @(0.107794, 0.323382, 0.06105731323805295, 0.06105731323805295, 3.0, 1.0, 2.758466666666667, 8.275399999999999, 0.1473530807662759, 0.1473530807662759, 0.2973461558666667, 2.6761154028, 0.008996983208938623, 0.008996983208938623)
    13 = @(0.107794, 0.323382, 0.06105731323805295, 0.06105731323805295, 3.0, 1.0, 0.1064366666666667, 0.31931, 0.005685684343896315, 0.005685684343896314, 0.01147323404666667, 0.10325910642, 0.0003471526099579708, 0.0003471526099579708)
}

foreach ($row in $newValues.Keys) {
    $values = $newValues[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range($columns[$i] + $row).Value = $values[$i]
    }
}

# Drop the obsolete Resolving-Mac sending-cluster rows (old rows 14-17) entirely.
$ws.Rows("14:17").Delete()
